# Update the missile policy parameters workbook with a working simulation:
# insert a new "Offensive Missile Success Probability" column in front of the
# existing "Defensive Missile Success Probability" column on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at H - this pushes the existing H column (Defensive
# Missile Success Probability, with its number format/style) to column I.
$ws.Columns.Item(8).Insert()

# New column header + values (offensive missile success probability).
$ws.Range("H1").Value = "Offensive Missile Success Probability"
$ws.Range("H2").Value = 0.9
$ws.Range("H3").Value = 0.9
$ws.Range("H2:H3").NumberFormat = "#\ ???/???"

# Resize the columns to fit their new contents.
$ws.Columns.Item(8).ColumnWidth = 30.833333333333332
$ws.Columns.Item(9).ColumnWidth = 24.5

# Move the selection, matching where the author left off editing.
$ws.Range("H4").Select()
